# "In Class Demonstrations" table update:
#   - adds two new regressor columns (U, $\pi$) next to the existing C/A, FFR
#     columns, turning the 2-column coefficient table into a 4-column one
#   - replaces the old "Constant"/"r2_adj" label rows with "U Lag"/"$\pi$ Lag"
#     rows of new coefficients
#   - every data cell is stored as text (even the numeric-looking ones), so
#     existing numeric cells need to be re-entered as text too
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grab the existing bold/bordered/centered header style (style index 1, already
# used by B1/C1 and the A-column labels) and stamp it onto the two new header
# cells before anything else touches them.
$ws.Range("B1").Copy()
$ws.Range("D1:E1").PasteSpecial(-4122)

# Wipe the old 3x5 table clean (keeps cell formatting in place).
$ws.Range("A1:E5").ClearContents()

# Force the whole data block to store values as text so strings that look
# like numbers (e.g. "0.055") are kept as shared-string text, not numbers.
$ws.Range("B2:E5").NumberFormat = "@"

# Header row.
$ws.Range("B1").Value = "C/A"
$ws.Range("C1").Value = "FFR"
$ws.Range("D1").Value = "U"
$ws.Range("E1").Value = '$\pi$'

# Row labels (column A).
$ws.Range("A2").Value = "C/A Lag"
$ws.Range("A3").Value = "FFR Lag"
$ws.Range("A4").Value = "U Lag"
$ws.Range("A5").Value = '$\pi$ Lag'

# C/A column coefficients.
$ws.Range("B2").Value = "-0.778***"
$ws.Range("B3").Value = "9.555***"
$ws.Range("B4").Value = "-0.286"
$ws.Range("B5").Value = "2.606***"

# FFR column coefficients.
$ws.Range("C2").Value = "-0.03***"
$ws.Range("C3").Value = "0.407***"
$ws.Range("C4").Value = "0.055"
$ws.Range("C5").Value = "0.044"

# U column coefficients.
$ws.Range("D2").Value = "0.036***"
$ws.Range("D3").Value = "-1.289***"
$ws.Range("D4").Value = "-0.135**"
$ws.Range("D5").Value = "0.037"

# $\pi$ column coefficients.
$ws.Range("E2").Value = "-0.023***"
$ws.Range("E3").Value = "0.392***"
$ws.Range("E4").Value = "-0.123**"
$ws.Range("E5").Value = "-0.669***"

# Drop the temporary text number-format again now that every value has been
# entered as text, so the data cells end up back at the default style.
$ws.Range("B2:E5").Style = "Normal"
